$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Cells.Item(7, 1).Value = 131092646
$ws.Cells.Item(7, 2).Value = 79244
$ws.Cells.Item(7, 5).Value = 6425
$ws.Cells.Item(7, 6).Value = 'Garnlav'
$ws.Cells.Item(7, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(7, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 17).Value = 585082
$ws.Cells.Item(7, 18).Value = 7060264
$ws.Cells.Item(7, 26).ClearContents()
$ws.Cells.Item(7, 28).ClearContents()
$ws.Cells.Item(7, 29).ClearContents()
$ws.Cells.Item(8, 1).Value = 131086958
$ws.Cells.Item(8, 17).Value = 585165
$ws.Cells.Item(8, 18).Value = 7060565
$ws.Cells.Item(8, 26).Value = '12:22'
$ws.Cells.Item(8, 28).Value = '12:22'
$ws.Cells.Item(9, 1).Value = 131085805
$ws.Cells.Item(9, 17).Value = 585215
$ws.Cells.Item(9, 18).Value = 7060513
$ws.Cells.Item(9, 19).Value = 10
$ws.Cells.Item(9, 26).Value = '12:01'
$ws.Cells.Item(9, 28).Value = '12:01'
$ws.Cells.Item(9, 49).Value = 'Kim Hultgren'
$ws.Cells.Item(9, 50).Value = 'Kim Hultgren'
$ws.Cells.Item(10, 1).Value = 131085613
$ws.Cells.Item(10, 2).Value = 57884
$ws.Cells.Item(10, 5).Value = 100109
$ws.Cells.Item(10, 6).Value = 'Tretåig hackspett'
$ws.Cells.Item(10, 7).Value = 'Picoides tridactylus'
$ws.Cells.Item(10, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(10, 13).Value = 'färska spår'
$ws.Cells.Item(10, 17).Value = 585222
$ws.Cells.Item(10, 18).Value = 7060481
$ws.Cells.Item(10, 19).Value = 15
$ws.Cells.Item(10, 26).Value = '11:52'
$ws.Cells.Item(10, 28).Value = '11:52'
$ws.Cells.Item(10, 29).Value = 'Färska ringhack gran'
$ws.Cells.Item(10, 49).Value = 'Daniel Rutschman'
$ws.Cells.Item(10, 50).Value = 'Daniel Rutschman'
$ws.Cells.Item(18, 1).Value = 131085737
$ws.Cells.Item(18, 2).Value = 79244
$ws.Cells.Item(18, 5).Value = 6425
$ws.Cells.Item(18, 6).Value = 'Garnlav'
$ws.Cells.Item(18, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(18, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(18, 17).Value = 585170
$ws.Cells.Item(18, 18).Value = 7060469
$ws.Cells.Item(18, 19).Value = 15
$ws.Cells.Item(18, 26).Value = '11:58'
$ws.Cells.Item(18, 28).Value = '11:58'
$ws.Cells.Item(18, 49).Value = 'Daniel Rutschman'
$ws.Cells.Item(18, 50).Value = 'Daniel Rutschman'
$ws.Cells.Item(19, 1).Value = 131085446
$ws.Cells.Item(19, 17).Value = 585301
$ws.Cells.Item(19, 18).Value = 7060488
$ws.Cells.Item(19, 19).Value = 10
$ws.Cells.Item(19, 26).Value = '11:41'
$ws.Cells.Item(19, 28).Value = '11:41'
$ws.Cells.Item(19, 49).Value = 'Kim Hultgren'
$ws.Cells.Item(19, 50).Value = 'Kim Hultgren'
$ws.Cells.Item(20, 1).Value = 131092554
$ws.Cells.Item(20, 2).Value = 57884
$ws.Cells.Item(20, 5).Value = 100109
$ws.Cells.Item(20, 6).Value = 'Tretåig hackspett'
$ws.Cells.Item(20, 7).Value = 'Picoides tridactylus'
$ws.Cells.Item(20, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(20, 13).Value = 'äldre spår'
$ws.Cells.Item(20, 17).Value = 585147
$ws.Cells.Item(20, 18).Value = 7060312
$ws.Cells.Item(20, 19).Value = 15
$ws.Cells.Item(20, 26).ClearContents()
$ws.Cells.Item(20, 28).ClearContents()
$ws.Cells.Item(20, 29).Value = 'Äldre ringhack, gran'
$ws.Cells.Item(20, 49).Value = 'Daniel Rutschman'
$ws.Cells.Item(20, 50).Value = 'Daniel Rutschman'
$ws.Cells.Item(21, 1).Value = 131092560
$ws.Cells.Item(21, 2).Value = 91805
$ws.Cells.Item(21, 5).Value = 1108
$ws.Cells.Item(21, 6).Value = 'Harticka'
$ws.Cells.Item(21, 7).Value = 'Pelloporus leporinus'
$ws.Cells.Item(21, 8).Value = '(Fr.) Krieglst.'
$ws.Cells.Item(21, 13).ClearContents()
$ws.Cells.Item(21, 17).Value = 585129
$ws.Cells.Item(21, 18).Value = 7060254
$ws.Cells.Item(21, 19).Value = 10
$ws.Cells.Item(21, 26).Value = '15:17'
$ws.Cells.Item(21, 28).Value = '15:17'
$ws.Cells.Item(21, 29).ClearContents()
$ws.Cells.Item(21, 49).Value = 'Kim Hultgren'
$ws.Cells.Item(21, 50).Value = 'Kim Hultgren'
$ws.Cells.Item(29, 1).Value = 131085240
$ws.Cells.Item(29, 2).Value = 57884
$ws.Cells.Item(29, 5).Value = 100109
$ws.Cells.Item(29, 6).Value = 'Tretåig hackspett'
$ws.Cells.Item(29, 7).Value = 'Picoides tridactylus'
$ws.Cells.Item(29, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(29, 13).Value = 'färska spår'
$ws.Cells.Item(29, 17).Value = 585289
$ws.Cells.Item(29, 18).Value = 7060293
$ws.Cells.Item(29, 26).Value = '11:16'
$ws.Cells.Item(29, 28).Value = '11:16'
$ws.Cells.Item(29, 29).Value = 'Ringhack på tall'
$ws.Cells.Item(30, 1).Value = 131085178
$ws.Cells.Item(30, 2).Value = 91829
$ws.Cells.Item(30, 5).Value = 5432
$ws.Cells.Item(30, 6).Value = 'Granticka'
$ws.Cells.Item(30, 7).Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Cells.Item(30, 8).Value = ''
$ws.Cells.Item(30, 13).ClearContents()
$ws.Cells.Item(30, 17).Value = 585225
$ws.Cells.Item(30, 18).Value = 7060258
$ws.Cells.Item(30, 19).Value = 10
$ws.Cells.Item(30, 26).Value = '11:08'
$ws.Cells.Item(30, 28).Value = '11:08'
$ws.Cells.Item(30, 29).ClearContents()
$ws.Cells.Item(30, 49).Value = 'Kim Hultgren'
$ws.Cells.Item(30, 50).Value = 'Kim Hultgren'
$ws.Cells.Item(31, 1).Value = 131085484
$ws.Cells.Item(31, 17).Value = 585303
$ws.Cells.Item(31, 18).Value = 7060488
$ws.Cells.Item(31, 19).Value = 15
$ws.Cells.Item(31, 26).ClearContents()
$ws.Cells.Item(31, 28).ClearContents()
$ws.Cells.Item(31, 29).Value = 'Färska ringhack gran'
$ws.Cells.Item(31, 49).Value = 'Daniel Rutschman'
$ws.Cells.Item(31, 50).Value = 'Daniel Rutschman'
$ws.Cells.Item(32, 1).Value = 131085569
$ws.Cells.Item(32, 2).Value = 79244
$ws.Cells.Item(32, 5).Value = 6425
$ws.Cells.Item(32, 6).Value = 'Garnlav'
$ws.Cells.Item(32, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(32, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(32, 17).Value = 585249
$ws.Cells.Item(32, 18).Value = 7060505
$ws.Cells.Item(33, 1).Value = 131087388
$ws.Cells.Item(33, 17).Value = 585131
$ws.Cells.Item(33, 18).Value = 7060627
$ws.Cells.Item(34, 1).Value = 131092590
$ws.Cells.Item(34, 17).Value = 585145
$ws.Cells.Item(34, 18).Value = 7060230
$ws.Cells.Item(34, 19).Value = 10
$ws.Cells.Item(34, 26).Value = '15:20'
$ws.Cells.Item(34, 28).Value = '15:20'
$ws.Cells.Item(34, 49).Value = 'Kim Hultgren'
$ws.Cells.Item(34, 50).Value = 'Kim Hultgren'
$ws.Cells.Item(35, 1).Value = 131092585
$ws.Cells.Item(35, 2).Value = 91805
$ws.Cells.Item(35, 5).Value = 1108
$ws.Cells.Item(35, 6).Value = 'Harticka'
$ws.Cells.Item(35, 7).Value = 'Pelloporus leporinus'
$ws.Cells.Item(35, 8).Value = '(Fr.) Krieglst.'
$ws.Cells.Item(35, 17).Value = 585130
$ws.Cells.Item(35, 18).Value = 7060263
$ws.Cells.Item(35, 19).Value = 15
$ws.Cells.Item(35, 26).ClearContents()
$ws.Cells.Item(35, 28).ClearContents()
$ws.Cells.Item(35, 49).Value = 'Daniel Rutschman'
$ws.Cells.Item(35, 50).Value = 'Daniel Rutschman'
